$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before the current row 148 ("Banquete" row of the
# 2021-10-25 (44494) Esparragos sample). This shifts the old rows 148-174 down to
# 151-177, matching the new dimension A1:R177.
$ws.Range("A148:A150").EntireRow.Insert()

# The rows that used to be 148/149/150 are now at 151/152/153 - duplicate them
# into the freshly inserted 148/149/150 rows so all the static columns
# (Mercado, Region, Codreg, Categoria, Variedad, Calidad, Unidad, Origen, Kg o
# Unidades, Clasificacion) come along, then overwrite the columns that changed
# for the new weekly sample (Fecha, Volumen, Precio minimo/maximo/promedio,
# Precio $/Kg).
$ws.Range("A148:R148").Value2 = $ws.Range("A151:R151").Value2
$ws.Range("A149:R149").Value2 = $ws.Range("A152:R152").Value2
$ws.Range("A150:R150").Value2 = $ws.Range("A153:R153").Value2

# Row 148 - Banquete
$ws.Cells.Item(148, 4).Value2 = 45204
$ws.Cells.Item(148, 10).Value2 = 740
$ws.Cells.Item(148, 11).Value2 = 1500
$ws.Cells.Item(148, 12).Value2 = 1600
$ws.Cells.Item(148, 13).Value2 = 1543
$ws.Cells.Item(148, 16).Value2 = 1543

# Row 149 - Primera
$ws.Cells.Item(149, 4).Value2 = 45204
$ws.Cells.Item(149, 10).Value2 = 700
$ws.Cells.Item(149, 11).Value2 = 1300
$ws.Cells.Item(149, 12).Value2 = 1400
$ws.Cells.Item(149, 13).Value2 = 1364
$ws.Cells.Item(149, 16).Value2 = 1364

# Row 150 - Segunda
$ws.Cells.Item(150, 4).Value2 = 45204
$ws.Cells.Item(150, 10).Value2 = 250
$ws.Cells.Item(150, 11).Value2 = 900
$ws.Cells.Item(150, 12).Value2 = 1000
$ws.Cells.Item(150, 13).Value2 = 960
$ws.Cells.Item(150, 16).Value2 = 960
